$wb = $excel.ActiveWorkbook

# Add the new "logical functions" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "logical functions"

# Populate the new sheet's content (mirrors the sibling "testOr"-style sheets).
$newSheet.Range("C7").Value = "Method boolean checkOr()"
$newSheet.Range("C8").Value = "return anyTrue(new boolean[]{true, false});"

# Match column C width used by the sheet (~34.57 Excel width units).
$newSheet.Columns.Item(3).ColumnWidth = 33.7083333

$newSheet.Range("C8").Select()
